# Apply 2023-04-02 crime data update across the workbook's per-neighborhood
# worksheets and the citywide/by-neighborhood summary sheets.
# Each (sheet, cell) pair below carries the 2023 (column J, and in a few
# cases 2021/H and 2022/I) full-year running totals forward by one day.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('J2').Value = 1618
$ws.Range('J3').Value = 1707
$ws.Range('H4').Value = 1689
$ws.Range('I4').Value = 1757
$ws.Range('J4').Value = 381
$ws.Range('J5').Value = 119
$ws.Range('J6').Value = 2223
$ws.Range('H7').Value = 26002
$ws.Range('I7').Value = 26203
$ws.Range('J7').Value = 6048

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('J2').Value = 25
$ws.Range('J7').Value = 72

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('J2').Value = 58
$ws.Range('J3').Value = 72
$ws.Range('J6').Value = 64
$ws.Range('J7').Value = 207

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('J3').Value = 26
$ws.Range('J6').Value = 21
$ws.Range('J7').Value = 76

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('J2').Value = 43
$ws.Range('J4').Value = 13
$ws.Range('J6').Value = 68
$ws.Range('J7').Value = 217

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('J3').Value = 10
$ws.Range('J7').Value = 41

$ws = $wb.Worksheets.Item('New City')
$ws.Range('J2').Value = 43
$ws.Range('J3').Value = 36
$ws.Range('J7').Value = 154

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('J2').Value = 50
$ws.Range('J8').Value = 371
$ws.Range('J11').Value = 77
$ws.Range('J15').Value = 79
$ws.Range('J19').Value = 209
$ws.Range('J20').Value = 125
$ws.Range('J29').Value = 334
$ws.Range('J31').Value = 41
$ws.Range('J33').Value = 259
$ws.Range('J36').Value = 93
$ws.Range('J37').Value = 207
$ws.Range('J42').Value = 234
$ws.Range('J44').Value = 49
$ws.Range('J47').Value = 54
$ws.Range('J48').Value = 50
$ws.Range('J50').Value = 34
$ws.Range('J51').Value = 78
$ws.Range('J52').Value = 139
$ws.Range('J54').Value = 118
$ws.Range('H63').Value = 238
$ws.Range('I63').Value = 195
$ws.Range('J63').Value = 26
$ws.Range('J65').Value = 154
$ws.Range('J66').Value = 16
$ws.Range('J67').Value = 217
$ws.Range('J69').Value = 16
$ws.Range('J76').Value = 94
$ws.Range('J77').Value = 44
$ws.Range('J79').Value = 185
$ws.Range('J83').Value = 148
$ws.Range('J85').Value = 276
$ws.Range('J90').Value = 66
$ws.Range('J91').Value = 69
$ws.Range('J93').Value = 27
$ws.Range('J96').Value = 72
$ws.Range('J99').Value = 76
$ws.Range('H101').Value = 26002
$ws.Range('I101').Value = 26203
$ws.Range('J101').Value = 6048

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('J2').Value = 48
$ws.Range('J3').Value = 48
$ws.Range('J7').Value = 148

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('J2').Value = 65
$ws.Range('J3').Value = 73
$ws.Range('J4').Value = 11
$ws.Range('J6').Value = 100
$ws.Range('J7').Value = 259

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('J2').Value = 33
$ws.Range('J3').Value = 21
$ws.Range('J6').Value = 57
$ws.Range('J7').Value = 118

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('J2').Value = 99
$ws.Range('J3').Value = 121
$ws.Range('J6').Value = 84
$ws.Range('J7').Value = 334

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('J2').Value = 51
$ws.Range('J3').Value = 59
$ws.Range('J5').Value = 14
$ws.Range('J7').Value = 209

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('J2').Value = 17
$ws.Range('J7').Value = 49

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('J3').Value = 7
$ws.Range('J4').Value = 6
$ws.Range('J7').Value = 50

$ws = $wb.Worksheets.Item('River North')
$ws.Range('J6').Value = 55
$ws.Range('J7').Value = 94

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('J2').Value = 69
$ws.Range('J3').Value = 109
$ws.Range('J6').Value = 75
$ws.Range('J7').Value = 276

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('J3').Value = 46
$ws.Range('J5').Value = 4
$ws.Range('J6').Value = 125
$ws.Range('J7').Value = 234

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range('J3').Value = 5
$ws.Range('J7').Value = 16

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('J4').Value = 5
$ws.Range('J7').Value = 69

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('J2').Value = 49
$ws.Range('J3').Value = 67
$ws.Range('J7').Value = 185

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('J4').Value = 15
$ws.Range('J6').Value = 36
$ws.Range('J7').Value = 125

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('J2').Value = 32
$ws.Range('J6').Value = 41
$ws.Range('J7').Value = 93

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range('J4').Value = 3
$ws.Range('J7').Value = 27

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('J2').Value = 33
$ws.Range('J6').Value = 53
$ws.Range('J7').Value = 139

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('J6').Value = 28
$ws.Range('J7').Value = 54

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('J3').Value = 20
$ws.Range('J6').Value = 40
$ws.Range('J7').Value = 79

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('J3').Value = 10
$ws.Range('J4').Value = 6
$ws.Range('J6').Value = 9
$ws.Range('J7').Value = 34

$ws = $wb.Worksheets.Item('North Center')
$ws.Range('J6').Value = 10
$ws.Range('J7').Value = 16

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('J2').Value = 27
$ws.Range('J7').Value = 77

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('J3').Value = 14
$ws.Range('J7').Value = 50

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('J2').Value = 121
$ws.Range('J3').Value = 127
$ws.Range('J7').Value = 371

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('J3').Value = 17
$ws.Range('J7').Value = 66

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('J3').Value = 23
$ws.Range('J7').Value = 78

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('J2').Value = 12
$ws.Range('J7').Value = 44
